$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> IP, add new PPS sheet right after it ---
$ip = $wb.Worksheets.Item(1)
$ip.Name = "IP"

$pps = $wb.Worksheets.Add($null, $ip)
$pps.Name = "PPS"

# ================= IP sheet =================
$ip.Select()

# New header cell E1 (match the bold/shaded header style used by A1:D1)
$ip.Range("E1").Value = "connection"
$ip.Range("D1").Copy()
$ip.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New rows 22-23 (nagios host)
$ip.Range("A22").Value = "nagios"
$ip.Range("B22").Value = "em0"
$ip.Range("C22").Value = "OPT1"
$ip.Range("D22").Value = "nat"
$ip.Range("E22").Value = "http://127.0.0.1:18033/nagios/"

$ip.Range("A23").Value = "nagios"
$ip.Range("B23").Value = "em1"
$ip.Range("C23").Value = "SOC"
$ip.Range("D23").Value = "10.10.3.3"

# Fixed typo 10.10.10.11 -> 10.10.11.1 (DEVWKS row, D21)
$ip.Range("D21").Value = "10.10.11.1"

# Column widths: keep col A + B:D as before, make col E wider
$ip.Columns.Item(5).ColumnWidth = 24.8

# Selection per diff
$ip.Range("A3").Select()

# ================= PPS sheet =================
$pps.Select()

$pps.Range("A1").Value = "host"
$pps.Range("B1").Value = "ipaddr"
$pps.Range("C1").Value = "protocol"
$pps.Range("D1").Value = "guest port"
$pps.Range("E1").Value = "host port"
$pps.Range("F1").Value = "service"

# Apply the same bold/shaded header style used on the IP sheet
$ip.Range("A1").Copy()
$pps.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ppsRows = @(
    @("rtr-ext", "nat", "https", 443, 55433, "pfsense web"),
    @("rtr-dmz", "nat", "https", 443, 55434, "pfsense web"),
    @("rtr-int", "nat", "https", 443, 55435, "pfsense web"),
    @("rtr-ent", "nat", "https", 443, 55436, "pfsense web"),
    @("rtr-dev", "nat", "https", 443, 55437, "pfsense web"),
    @("nagios",  "nat", "http",  80,  18033, "nagios web")
)

$r = 2
foreach ($row in $ppsRows) {
    $pps.Cells.Item($r, 1).Value = $row[0]
    $pps.Cells.Item($r, 2).Value = $row[1]
    $pps.Cells.Item($r, 3).Value = $row[2]
    $pps.Cells.Item($r, 4).Value = $row[3]
    $pps.Cells.Item($r, 5).Value = $row[4]
    $pps.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

$pps.Range("B24").Select()

# Restore the IP sheet as the active sheet/tab
$ip.Select()
